$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# Row 7: DMA-80 THG Data_Cleandata -> B7 = "Yes" (clear the yellow "to fill" formatting entirely)
$ws.Range("B7").ClearFormats()
$ws.Range("B7").Value = "Yes"

# Row 9: L.Robbins Lab Data 2021-2023 ... -> B9 = "Yes" (fill cleared to "no fill")
$ws.Range("B9").ClearFormats()
$ws.Range("B9").Interior.ColorIndex = -4142
$ws.Range("B9").Value = "Yes"

# Row 12: MeHg_2021_ID_Samples_results -> B12 = "No", C12 = note (fill cleared to "no fill")
$ws.Range("B12").ClearFormats()
$ws.Range("B12").Interior.ColorIndex = -4142
$ws.Range("B12").Value = "No"
$ws.Range("C12").Value = "Was already included in updated PPR_2006_to_2024 Master File"

# Row 23: SDWS 2015-2022 Baulch -> B23 = "Yes", C23 = new note (clear formatting entirely)
$ws.Range("B23").ClearFormats()
$ws.Range("B23").Value = "Yes"
$ws.Range("C23").Value = "Merged with SDWS Long Term Monitoring data"

# Row 28: TICTOC Results May 30, 2023 -> B28 = "Yes" (fill cleared to "no fill")
$ws.Range("B28").ClearFormats()
$ws.Range("B28").Interior.ColorIndex = -4142
$ws.Range("B28").Value = "Yes"

# Update the active selection to B6, matching the saved view state
$ws.Range("B6").Select()
